$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "obj1"
$ws.Range("A2").Value = "obj2"
$ws.Range("A3").Value = "obj3"

$ws.Range("L15").Select()
